$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '30.323.90'
Set-TextValue $ws.Range("E2") '  -0.09%  '
Set-TextValue $ws.Range("D3") '1.931.28'
Set-TextValue $ws.Range("E3") '  -0.53%  '
Set-TextValue $ws.Range("E4") '  +0.05%  '
Set-TextValue $ws.Range("D5") '0.7379'
Set-TextValue $ws.Range("E5") '  +1.99%  '
Set-TextValue $ws.Range("D6") '249.98'
Set-TextValue $ws.Range("E6") '  -0.73%  '
Set-TextValue $ws.Range("D7") '1.001'
Set-TextValue $ws.Range("E7") '  +0.05%  '
Set-TextValue $ws.Range("D8") '0.3218'
Set-TextValue $ws.Range("E8") '  -3.92%  '
Set-TextValue $ws.Range("D9") '27.79'
Set-TextValue $ws.Range("E9") '  -3.96%  '
Set-TextValue $ws.Range("D10") '0.07102'
Set-TextValue $ws.Range("E10") '  -4.02%  '
Set-TextValue $ws.Range("D11") '0.7878'
Set-TextValue $ws.Range("E11") '  -3.87%  '
Set-TextValue $ws.Range("D12") '0.08023'
Set-TextValue $ws.Range("E12") '  -1.45%  '
Set-TextValue $ws.Range("D13") '1.933.58'
Set-TextValue $ws.Range("E13") '  -0.47%  '
Set-TextValue $ws.Range("D14") '5.390'
Set-TextValue $ws.Range("E14") '  -1.92%  '
Set-TextValue $ws.Range("D15") '94.68'
Set-TextValue $ws.Range("E15") '  -0.62%  '
Set-TextValue $ws.Range("D16") '14.55'
Set-TextValue $ws.Range("E16") '  -2.66%  '
Set-TextValue $ws.Range("D17") '30.333.14'
Set-TextValue $ws.Range("D18") '253.82'
Set-TextValue $ws.Range("E18") '  +0.03%  '
Set-TextValue $ws.Range("D19") '0.000008064'
Set-TextValue $ws.Range("E19") '  -4.10%  '
Set-TextValue $ws.Range("E20") '  -2.64%  '
Set-TextValue $ws.Range("D21") '2.186.72'
Set-TextValue $ws.Range("E21") '  -0.54%  '
Set-TextValue $ws.Range("E22") '  +0.04%  '
Set-TextValue $ws.Range("D23") '1.001'
Set-TextValue $ws.Range("E23") '  -0.08%  '
Set-TextValue $ws.Range("D24") '6.829'
Set-TextValue $ws.Range("E24") '  -2.33%  '
Set-TextValue $ws.Range("D25") '9.565'
Set-TextValue $ws.Range("E25") '  -3.07%  '
Set-TextValue $ws.Range("D26") '164.04'
Set-TextValue $ws.Range("E26") '  +0.47%  '
Set-TextValue $ws.Range("D27") '19.08'
Set-TextValue $ws.Range("E27") '  -1.54%  '
Set-TextValue $ws.Range("D28") '2.282'
Set-TextValue $ws.Range("E28") '  -5.63%  '
Set-TextValue $ws.Range("D29") '0.1319'
Set-TextValue $ws.Range("E29") '  -0.68%  '
Set-TextValue $ws.Range("D30") '1.353'
Set-TextValue $ws.Range("E30") '  +0.55%  '
Set-TextValue $ws.Range("D31") '1.534'
Set-TextValue $ws.Range("E31") '  -2.76%  '
Set-TextValue $ws.Range("D32") '4.416'
Set-TextValue $ws.Range("E32") '  -1.15%  '
Set-TextValue $ws.Range("D33") '4.149'
Set-TextValue $ws.Range("E33") '  -2.65%  '
Set-TextValue $ws.Range("D34") '0.05120'
Set-TextValue $ws.Range("E34") '  -3.28%  '
Set-TextValue $ws.Range("D35") '1.289'
Set-TextValue $ws.Range("E35") '  -1.66%  '
Set-TextValue $ws.Range("D36") '0.7467'
Set-TextValue $ws.Range("E36") '  -1.56%  '
Set-TextValue $ws.Range("D37") '2.769'
Set-TextValue $ws.Range("E37") '  +0.93%  '
Set-TextValue $ws.Range("D38") '0.01977'
Set-TextValue $ws.Range("E38") '  -0.88%  '
Set-TextValue $ws.Range("E39") '  -1.68%  '
Set-TextValue $ws.Range("D40") '77.94'
Set-TextValue $ws.Range("E40") '  -4.02%  '
Set-TextValue $ws.Range("D41") '6.402'
Set-TextValue $ws.Range("E41") '  -3.32%  '
Set-TextValue $ws.Range("D42") '0.4504'
Set-TextValue $ws.Range("E42") '  -1.76%  '
Set-TextValue $ws.Range("D43") '1.988'
Set-TextValue $ws.Range("E43") '  -2.96%  '
Set-TextValue $ws.Range("D44") '0.8462'
Set-TextValue $ws.Range("E44") '  -0.14%  '
Set-TextValue $ws.Range("D45") '1.001'
Set-TextValue $ws.Range("E45") '  +0.03%  '
Set-TextValue $ws.Range("D46") '101.29'
Set-TextValue $ws.Range("E46") '  -1.48%  '
Set-TextValue $ws.Range("D47") '7.531'
Set-TextValue $ws.Range("E47") '  -0.05%  '
Set-TextValue $ws.Range("D48") '9.769'
Set-TextValue $ws.Range("E48") '  -0.50%  '
Set-TextValue $ws.Range("D49") '982.27'
Set-TextValue $ws.Range("E49") '  +11.52%  '
Set-TextValue $ws.Range("D50") '37.02'
Set-TextValue $ws.Range("E50") '  +0.17%  '
Set-TextValue $ws.Range("D51") '0.06077'
Set-TextValue $ws.Range("E51") '  +0.46%  '

Write-Host "Applied cryptos update"